$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Widened item names to 11": four item-size offsets get their low nibble
# bumped from 8 to a (e.g. 0818 -> 081a), widening those entries.
# A leading apostrophe forces the values to stay text (matching the
# existing t="s" shared-string cell type/style) instead of being
# auto-converted to numbers.
$ws.Range("D10").Value = "'081a"
$ws.Range("D12").Value = "'041a"
$ws.Range("D15").Value = "'0319"
$ws.Range("D20").Value = "'151a"

$excel.CalculateFullRebuild()

# Update the view state left over from editing: the selection ends up on I19.
$ws.Range("I19").Select()
